# Cronograma.xlsx update — "Validacao Usuario, competencia, situacao"
#
# Updates progress figures for the "View / Cadastro" column (Competencia,
# Hidrometro UC, Situacao, Usuario rows) on the Plan1 schedule sheet, and
# clears the now-obsolete "Modificar a estrutura de validacao para MD5"
# note from the Melhorias sheet now that the validation item is done.

$wb = $excel.ActiveWorkbook

$plan1 = $wb.Worksheets.Item("Plan1")
$melhorias = $wb.Worksheets.Item("Melhorias")

# --- Plan1: bump "Cadastro" progress (column G) for the rows touched by
#     this validation pass, and mark Hidrometro UC / Modelo as started. ---

# Competencia (row 7)
$plan1.Range("G7").Value = 0.95

# Hidrometro UC (row 11) - Modelo column now started
$plan1.Range("B11").Value = 0

# Situacao (row 19)
$plan1.Range("G19").Value = 0.95

# Usuario (row 21)
$plan1.Range("G21").Value = 0.95

# --- Melhorias: the MD5 validation note is resolved, drop it ---
$melhorias.Range("B4").ClearContents()

# --- Leave the workbook focused on Plan1, matching the user's last view ---
$melhorias.Activate()
$null = $melhorias.Range("B5").Select()

$plan1.Activate()
$null = $plan1.Range("G22").Select()
